# Apply crypto price/volume updates (GitHub Actions data refresh)
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").NumberFormat = "@"
$ws.Range("D2").Value = "34.499.95"
$ws.Range("E2").Value = "  -0.31%  "
$ws.Range("D3").NumberFormat = "@"
$ws.Range("D3").Value = "1.807.24"
$ws.Range("E3").Value = "  -0.61%  "
$ws.Range("E4").Value = "  +0.19%  "
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = "227.90"
$ws.Range("E5").Value = "  -0.19%  "
$ws.Range("D6").NumberFormat = "@"
$ws.Range("D6").Value = "0.580"
$ws.Range("E6").Value = "  +3.80%  "
$ws.Range("E7").Value = "  +0.17%  "
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = "36.80"
$ws.Range("E8").Value = "  +5.91%  "
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = "0.299"
$ws.Range("E9").Value = "  -0.62%  "
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = "0.0693"
$ws.Range("E10").Value = "  -0.37%  "
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = "0.0964"
$ws.Range("E11").Value = "  +1.39%  "
$ws.Range("D12").NumberFormat = "@"
$ws.Range("D12").Value = "2.068.04"
$ws.Range("E12").Value = "  -0.54%  "
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = "11.60"
$ws.Range("E13").Value = "  +2.04%  "
$ws.Range("B14").Value = "Polygon"
$ws.Range("C14").Value = "https://coinranking.com/coin/uW2tk-ILY0ii+polygon-matic"
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = "0.655"
$ws.Range("E14").Value = "  +1.46%  "
$ws.Range("B15").Value = "WrappedEther"
$ws.Range("C15").Value = "https://coinranking.com/coin/Mtfb0obXVh59u+wrappedether-weth"
$ws.Range("D15").NumberFormat = "@"
$ws.Range("D15").Value = "1.802.61"
$ws.Range("E15").Value = "  -0.85%  "
$ws.Range("E16").Value = "  +3.27%  "
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = "34.494.63"
$ws.Range("E17").Value = "  -0.31%  "
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = "69.62"
$ws.Range("E18").Value = "  +0.69%  "
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = "245.38"
$ws.Range("E19").Value = "  -0.84%  "
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = "0.0₃0791"
$ws.Range("E20").Value = "  -1.36%  "
$ws.Range("D21").NumberFormat = "@"
$ws.Range("D21").Value = "11.67"
$ws.Range("E21").Value = "  +0.88%  "
$ws.Range("E22").Value = "  +0.15%  "
$ws.Range("E23").Value = "  -0.52%  "
$ws.Range("E24").Value = "  +5.24%  "
$ws.Range("D25").NumberFormat = "@"
$ws.Range("D25").Value = "172.43"
$ws.Range("E25").Value = "  -0.18%  "
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = "7.99"
$ws.Range("E26").Value = "  +7.11%  "
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = "16.92"
$ws.Range("E27").Value = "  +0.94%  "
$ws.Range("E28").Value = "  +1.59%  "
$ws.Range("E29").Value = "  +0.07%  "
$ws.Range("E30").Value = "  -0.50%  "
$ws.Range("B31").Value = "Hedera"
$ws.Range("C31").Value = "https://coinranking.com/coin/jad286TjB+hedera-hbar"
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = "0.0530"
$ws.Range("E31").Value = "  -0.21%  "
$ws.Range("B32").Value = "Filecoin"
$ws.Range("C32").Value = "https://coinranking.com/coin/ymQub4fuB+filecoin-fil"
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = "3.85"
$ws.Range("E32").Value = "  +0.08%  "
$ws.Range("E33").Value = "  -0.12%  "
$ws.Range("E34").Value = "  -1.84%  "
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = "1.395.83"
$ws.Range("E35").Value = "  -1.73%  "
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = "0.673"
$ws.Range("E36").Value = "  -0.60%  "
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = "2.46"
$ws.Range("E37").Value = "  -5.15%  "
$ws.Range("E38").Value = "  -0.58%  "
$ws.Range("E39").Value = "  -0.31%  "
$ws.Range("B40").Value = "ARBITRUM"
$ws.Range("C40").Value = "https://coinranking.com/coin/1Uo6s62Oc+arbitrum-arb"
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = "0.971"
$ws.Range("E40").Value = "  +1.71%  "
$ws.Range("B41").Value = "Aave"
$ws.Range("C41").Value = "https://coinranking.com/coin/ixgUfzmLR+aave-aave"
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = "83.18"
$ws.Range("E41").Value = "  -3.18%  "
$ws.Range("E42").Value = "  -0.68%  "
$ws.Range("E43").Value = "  +0.87%  "
$ws.Range("E44").Value = "  +7.74%  "
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = "13.50"
$ws.Range("E45").Value = "  -2.22%  "
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = "6.04"
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = "0.0501"
$ws.Range("E47").Value = "  -5.05%  "
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = "1.969.66"
$ws.Range("E48").Value = "  -0.54%  "
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = "104.46"
$ws.Range("E49").Value = "  -1.22%  "
$ws.Range("E50").Value = "  +0.11%  "
$ws.Range("E51").Value = "  -3.06%  "
